$d = $word.ActiveDocument

# 1) Update the "Choose vehicle type" line from In Progress to Completed.
$d.Content.Find.Execute(
    "[~] (In Progress) Choose vehicle type (car/okada/pragya/aboboyaa)",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "[x] (Completed) Choose vehicle type (car/okada/pragya/aboboyaa)", 2)

# 2) Insert four new checklist items right after the "Push branch ..." line,
#    before the blank paragraph that precedes "Artifacts created:".
$idx = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $t = $d.Paragraphs($i).Range.Text
    if ($t -like "*Push branch feature/rider-backend-wiring from Windows (VM DNS blocked)*") {
        $idx = $i
        break
    }
}

$lines = @(
    "[x] Rider profile updates call backend (/users/me)",
    "[x] App config supports API base URL + mock toggle via dart-define",
    "[x] Delivery request hits backend (/deliveries/request)",
    "[x] Home shows API base URL + backend reachability"
)

$insertAfterIdx = $idx
foreach ($line in $lines) {
    $d.Paragraphs($insertAfterIdx).Range.InsertParagraphAfter()
    $insertAfterIdx = $insertAfterIdx + 1
    $d.Paragraphs($insertAfterIdx).Range.Text = $line
}
